$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")

# Delete row 89 (Caso -500 / Castañares 5656) entirely; rows below shift up.
$ws.Rows.Item(89).Delete()
